$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (A1:F1) stays the same content (nis, nama, gambar, quote, jurusan, d_kelas)
$ws.Range("A1").Value = "nis"
$ws.Range("F1").Value = "d_kelas"

# Numeric columns A and F for data rows
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4

$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 4
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1

# Fill the string columns column-by-column (matches the order strings were
# (re)written so the shared-strings table ends up in the expected order)
$ws.Range("B1").Value = "nama"
$ws.Range("B2").Value = "not error"
$ws.Range("B3").Value = "not error"
$ws.Range("B4").Value = "not error"
$ws.Range("B5").Value = "not error"

$ws.Range("E1").Value = "jurusan"
$ws.Range("E2").Value = "tei"
$ws.Range("E3").Value = "tei"
$ws.Range("E4").Value = "TKJ"
$ws.Range("E5").Value = "TKJ"

$ws.Range("C1").Value = "gambar"
$ws.Range("C2").Value = "Si_A.jpg"
$ws.Range("C3").Value = "Si_B.jpg"
$ws.Range("C4").Value = "Si_C.jpg"
$ws.Range("C5").Value = "si_D.jpg"

$ws.Range("D1").Value = "quote"
$ws.Range("D2").Value = "ah masa"
$ws.Range("D3").Value = "ah masa"
$ws.Range("D4").Value = "ah masa"
$ws.Range("D5").Value = "ah masa"

# Update selection to match diff (activeCell E3)
$ws.Range("E3").Select()
